$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated "Price" (D) and "Volume(1h)" (E) figures scraped
# for this run. Several of the new D-column price strings parse as
# plain numbers (e.g. "1.000", "237.62", "4.600") -- set the cell to
# Text format first so Excel keeps the literal digits (trailing
# zeros, "NN.NNN.NN" thousand-dot grouping) instead of silently
# coercing them to the General number type and losing formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.861.92'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.738.52'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.62'
$ws.Range("E5").Value = '  +3.43%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5184'
$ws.Range("E7").Value = '  -1.67%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06150'
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.742.17'
$ws.Range("E10").Value = '  +0.33%  '
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.93'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6413'
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.600'
$ws.Range("E14").Value = '  +1.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.46'
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.887.63'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("E19").Value = '  +1.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006755'
$ws.Range("E20").Value = '  +1.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.966.80'
$ws.Range("E21").Value = '  +0.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.273'
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.618'
$ws.Range("E23").Value = '  -1.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.259'
$ws.Range("E24").Value = '  +1.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.97'
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.515'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.14'
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.757'
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '105.70'
$ws.Range("E29").Value = '  +3.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.904'
$ws.Range("E30").Value = '  +5.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08272'
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.684'
$ws.Range("E32").Value = '  +4.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04604'
$ws.Range("E33").Value = '  +2.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.649'
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("E35").Value = '  +1.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6163'
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.679'
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01603'
$ws.Range("E38").Value = '  +2.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.923'
$ws.Range("E39").Value = '  +1.34%  '
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.36'
$ws.Range("E41").Value = '  -2.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3835'
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("E43").Value = '  +2.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.979'
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.224'
$ws.Range("E46").Value = '  +0.34%  '
$ws.Range("E47").Value = '  -1.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.69'
$ws.Range("E48").Value = '  +2.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.39'
$ws.Range("E49").Value = '  +1.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.547'
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3404'
$ws.Range("E51").Value = '  +0.28%  '
